# Auto-generated edit script applying value updates per the OOXML diff.
# Each sheet is selected once; cell values are updated or cleared to match target state.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 35988.43
$ws.Range("J3").Value = 35988.43
$ws.Range("L3").Value = 35988.43
$ws.Range("N3").Value = -36216.43
$ws.Range("H13").Value = 4570.857
$ws.Range("I13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("M13").ClearContents()
$ws.Range("H17").Value = 732.7368
$ws.Range("J17").Value = 777.03845
$ws.Range("L17").Value = 2331.11535
$ws.Range("N17").Value = -2667.11535
$ws.Range("H28").Value = 1243.4286
$ws.Range("I28").Value = 1021
$ws.Range("K28").Value = 1021
$ws.Range("M28").Value = -536
$ws.Range("H31").Value = 588.8570999999999
$ws.Range("I31").Value = 588.8570999999999
$ws.Range("K31").Value = 1766.5713
$ws.Range("M31").Value = -1536.5713
$ws.Range("H33").Value = 2790.1516
$ws.Range("I33").Value = 2986.16
$ws.Range("K33").Value = 2986.16
$ws.Range("M33").Value = -2757.16
$ws.Range("H41").Value = 1262.4445
$ws.Range("I41").Value = 1375.8889
$ws.Range("J41").Value = 1149
$ws.Range("K41").Value = 1375.8889
$ws.Range("L41").Value = 1149
$ws.Range("M41").Value = -935.8888999999999
$ws.Range("N41").Value = -2029
$ws.Range("H68").Value = 80262.5
$ws.Range("J68").Value = 80262.5
$ws.Range("L68").Value = 80262.5
$ws.Range("N68").Value = -81760.5
$ws.Range("H71").Value = 80262.5
$ws.Range("J71").Value = 80262.5
$ws.Range("L71").Value = 240787.5
$ws.Range("N71").Value = -248275.5
$ws.Range("H87").Value = 124959.5
$ws.Range("J87").Value = 124959.5
$ws.Range("L87").Value = 124959.5
$ws.Range("N87").Value = -127455.5
$ws.Range("H90").Value = 124959.5
$ws.Range("J90").Value = 124959.5
$ws.Range("L90").Value = 374878.5
$ws.Range("N90").Value = -387358.5
$ws.Range("H93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("L93").ClearContents()
$ws.Range("N93").Value = 0
$ws.Range("H102").Value = 35988.43
$ws.Range("J102").Value = 35988.43
$ws.Range("L102").Value = 35988.43
$ws.Range("N102").Value = -42478.43
$ws.Range("H107").Value = 351.77777
$ws.Range("I107").Value = 270.75
$ws.Range("K107").Value = 270.75
$ws.Range("M107").Value = 1649.25
$ws.Range("H116").Value = 7315.273
$ws.Range("J116").Value = 8353.6
$ws.Range("L116").Value = 8353.6
$ws.Range("N116").Value = -15237.6
$ws.Range("H132").Value = 1457.8334
$ws.Range("I132").Value = 1378.303
$ws.Range("J132").Value = 2332.6667
$ws.Range("K132").Value = 4134.909000000001
$ws.Range("L132").Value = 6998.000100000001
$ws.Range("M132").Value = -1604.909000000001
$ws.Range("N132").Value = -12058.0001

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 11114456
$ws.Range("I32").Value = 11907351
$ws.Range("K32").Value = 11907351
$ws.Range("M32").Value = -11907064
$ws.Range("H61").Value = 9639293
$ws.Range("I61").Value = 11367625
$ws.Range("J61").Value = 133463.38
$ws.Range("K61").Value = 11367625
$ws.Range("L61").Value = 133463.38
$ws.Range("M61").Value = -11367413
$ws.Range("N61").Value = -133887.38
$ws.Range("H97").Value = 950.35297
$ws.Range("I97").Value = 950.35297
$ws.Range("K97").Value = 950.35297
$ws.Range("M97").Value = -454.35297
$ws.Range("H110").Value = 1323.8334
$ws.Range("I110").Value = 1323.8334
$ws.Range("K110").Value = 1323.8334
$ws.Range("M110").Value = 721.1666
$ws.Range("H132").Value = 4504.7383
$ws.Range("I132").Value = 1900.5428
$ws.Range("J132").Value = 17525.715
$ws.Range("K132").Value = 5701.6284
$ws.Range("L132").Value = 52577.145
$ws.Range("M132").Value = -3171.6284
$ws.Range("N132").Value = -57637.145
$ws.Range("H136").Value = 9639293
$ws.Range("I136").Value = 11367625
$ws.Range("J136").Value = 133463.38
$ws.Range("K136").Value = 34102875
$ws.Range("L136").Value = 400390.14
$ws.Range("M136").Value = -34100325
$ws.Range("N136").Value = -405490.14

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3121.15
$ws.Range("I20").Value = 3027.7334
$ws.Range("J20").Value = 3401.4
$ws.Range("K20").Value = 3027.7334
$ws.Range("L20").Value = 3401.4
$ws.Range("M20").Value = -2780.7334
$ws.Range("N20").Value = -3895.4
$ws.Range("H80").Value = 765.2
$ws.Range("J80").Value = 619.2727
$ws.Range("L80").Value = 619.2727
$ws.Range("N80").Value = -2615.2727
$ws.Range("H83").Value = 765.2
$ws.Range("J83").Value = 619.2727
$ws.Range("L83").Value = 3096.3635
$ws.Range("N83").Value = -13080.3635
$ws.Range("H99").Value = 2202.1333
$ws.Range("I99").Value = 1769.3334
$ws.Range("K99").Value = 1769.3334
$ws.Range("M99").Value = -271.3334
$ws.Range("H107").Value = 4314.727
$ws.Range("I107").Value = 3380.125
$ws.Range("K107").Value = 3380.125
$ws.Range("M107").Value = -1460.125

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 900
$ws.Range("I16").Value = 900
$ws.Range("K16").Value = 900
$ws.Range("M16").Value = -613
$ws.Range("H88").Value = 0
$ws.Range("J88").Value = 0
$ws.Range("L88").ClearContents()
$ws.Range("N88").Value = 0
$ws.Range("H91").Value = 0
$ws.Range("J91").Value = 0
$ws.Range("L91").ClearContents()
$ws.Range("N91").Value = 0
$ws.Range("H107").Value = 508.5
$ws.Range("J107").Value = 894.5
$ws.Range("L107").Value = 894.5
$ws.Range("N107").Value = -4734.5
$ws.Range("H113").Value = 900
$ws.Range("I113").Value = 900
$ws.Range("K113").Value = 900
$ws.Range("M113").Value = 1270
$ws.Range("H132").Value = 1382.8
$ws.Range("I132").Value = 1192.421
$ws.Range("J132").Value = 5000
$ws.Range("K132").Value = 3577.263
$ws.Range("L132").Value = 15000
$ws.Range("M132").Value = -1047.263
$ws.Range("N132").Value = -20060
$ws.Range("H134").Value = 421355.12
$ws.Range("I134").Value = 667445.7
$ws.Range("J134").Value = 11204.223
$ws.Range("K134").Value = 2002337.1
$ws.Range("L134").Value = 33612.669
$ws.Range("M134").Value = -1999802.1
$ws.Range("N134").Value = -38682.669

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 12657.857
$ws.Range("I3").Value = 7735.4
$ws.Range("J3").Value = 15392.556
$ws.Range("K3").Value = 23206.2
$ws.Range("L3").Value = 46177.66800000001
$ws.Range("M3").Value = -23094.2
$ws.Range("N3").Value = -46401.66800000001
$ws.Range("H37").Value = 0
$ws.Range("J37").Value = 0
$ws.Range("L37").ClearContents()
$ws.Range("N37").Value = 0
$ws.Range("H44").Value = 37863.734
$ws.Range("I44").Value = 37863.734
$ws.Range("K44").Value = 113591.202
$ws.Range("M44").Value = -113193.202
$ws.Range("H56").Value = 9999.666999999999
$ws.Range("I56").Value = 9999.666999999999
$ws.Range("K56").Value = 9999.666999999999
$ws.Range("M56").Value = -9469.666999999999
$ws.Range("H107").Value = 527.3333
$ws.Range("I107").Value = 437.1111
$ws.Range("J107").Value = 617.55554
$ws.Range("K107").Value = 1311.3333
$ws.Range("L107").Value = 1852.66662
$ws.Range("M107").Value = 608.6667
$ws.Range("N107").Value = -5692.66662
$ws.Range("H112").Value = 10657.143
$ws.Range("I112").Value = 9933.333000000001
$ws.Range("J112").Value = 15000
$ws.Range("K112").Value = 29799.999
$ws.Range("L112").Value = 45000
$ws.Range("M112").Value = -28691.999
$ws.Range("N112").Value = -47216
$ws.Range("H133").Value = 5587.1665
$ws.Range("I133").Value = 5116.222
$ws.Range("K133").Value = 15348.666
$ws.Range("M133").Value = -10288.666
$ws.Range("H140").Value = 117545.46
$ws.Range("I140").Value = 127059.414
$ws.Range("K140").Value = 381178.242
$ws.Range("M140").Value = -375998.242

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 2980.5
$ws.Range("I113").Value = 1011
$ws.Range("K113").Value = 1011
$ws.Range("M113").Value = 1159

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 4200.6665
$ws.Range("I22").Value = 6001
$ws.Range("J22").Value = 600
$ws.Range("K22").Value = 6001
$ws.Range("L22").Value = 600
$ws.Range("M22").Value = -5706
$ws.Range("N22").Value = -1190
$ws.Range("H27").Value = 4200.6665
$ws.Range("I27").Value = 6001
$ws.Range("J27").Value = 600
$ws.Range("K27").Value = 6001
$ws.Range("L27").Value = 600
$ws.Range("M27").Value = -5894
$ws.Range("N27").Value = -814
$ws.Range("H54").Value = 0
$ws.Range("J54").Value = 0
$ws.Range("L54").ClearContents()
$ws.Range("N54").Value = 0
$ws.Range("H100").Value = 1367.4
$ws.Range("I100").Value = 1367.4
$ws.Range("K100").Value = 1367.4
$ws.Range("M100").Value = -826.4000000000001
$ws.Range("H132").Value = 380951.6
$ws.Range("I132").Value = 358842.53
$ws.Range("J132").Value = 1000005
$ws.Range("K132").Value = 1076527.59
$ws.Range("L132").Value = 3000015
$ws.Range("M132").Value = -1073997.59
$ws.Range("N132").Value = -3005075

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 38495
$ws.Range("J54").Value = 38495
$ws.Range("L54").Value = 38495
$ws.Range("N54").Value = -39535
$ws.Range("H132").Value = 1616.0667
$ws.Range("I132").Value = 1399.92
$ws.Range("J132").Value = 2696.8
$ws.Range("K132").Value = 4199.76
$ws.Range("L132").Value = 8090.400000000001
$ws.Range("M132").Value = -1669.76
$ws.Range("N132").Value = -13150.4
